# Update Name of Algo
# Applies the numeric updates described in the diff to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3").Value = 6.079099999999993
$ws.Range("E4").Value = 13.1374
$ws.Range("E6").Value = 12.19219999999999
$ws.Range("A9").Value = -20.29509999999998
$ws.Range("E10").Value = 12.17529999999999
$ws.Range("B11").Value = 5.271999999999998
$ws.Range("E11").Value = 13.45229999999999
$ws.Range("A18").Value = -22.864
$ws.Range("A20").Value = -22.08860000000002
$ws.Range("D21").Value = -7.704400000000005
$ws.Range("E21").Value = 13.48300000000001
